# Update countries & provincias Spain
#
# The "Pais" sheet lists COVID-19 stats per country, one row per country,
# already sorted by "Casos totales" (column B) descending. This update
# applies a newer data snapshot: most changed rows simply get refreshed
# B:H figures, but a handful of countries' updated totals move them past
# their neighbours in the ranking, so those rows' country name (column A)
# changes too while the row number (and rank position) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Alemania) - refreshed figures, no reorder
$ws.Range("B8").Value = 135843
$ws.Range("C8").Value = 1090
$ws.Range("E8").Value = 54953
$ws.Range("G8").Value = 86
$ws.Range("H8").Value = 3890

# Row 14 -> Brasil
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 30425
$ws.Range("C14").Value = 1815
$ws.Range("D14").Value = 14026
$ws.Range("E14").Value = 14475
$ws.Range("F14").Value = 6634
$ws.Range("G14").Value = 167
$ws.Range("H14").Value = 1924

# Row 15 -> Canada
$ws.Range("A15").Value = "Canada"
$ws.Range("B15").Value = 29925
$ws.Range("C15").Value = 1546
$ws.Range("D15").Value = 9671
$ws.Range("E15").Value = 19063
$ws.Range("F15").Value = 557
$ws.Range("H15").Value = 1191

# Row 16 -> Paises Bajos
$ws.Range("A16").Value = "Paises Bajos"
$ws.Range("C16").Value = 1061
$ws.Range("D16").Value = 250
$ws.Range("E16").Value = 25649
$ws.Range("F16").Value = 1279
$ws.Range("G16").Value = 181
$ws.Range("H16").Value = 3315

# Row 20 (Austria) - refreshed figures, no reorder
$ws.Range("E20").Value = 5078
$ws.Range("F20").Value = 238
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 410

# Row 21 -> Irlanda
$ws.Range("A21").Value = "Irlanda"
$ws.Range("B21").Value = 13271
$ws.Range("C21").Value = 724
$ws.Range("D21").Value = 77
$ws.Range("E21").Value = 12708
$ws.Range("F21").Value = 158
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = 486

# Row 22 -> India
$ws.Range("A22").Value = "India"
$ws.Range("B22").Value = 12759
$ws.Range("C22").Value = 389
$ws.Range("D22").Value = 1514
$ws.Range("E22").Value = 10822
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 423

# Row 23 -> Israel
$ws.Range("A23").Value = "Israel"
$ws.Range("B23").Value = 12758
$ws.Range("C23").Value = 257
$ws.Range("D23").Value = 2818
$ws.Range("E23").Value = 9798
$ws.Range("F23").Value = 181
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 142

# Row 55 -> Sudafrica
$ws.Range("A55").Value = "Sudafrica"
$ws.Range("B55").Value = 2605
$ws.Range("C55").Value = 99
$ws.Range("D55").Value = 903
$ws.Range("E55").Value = 1654
$ws.Range("F55").Value = 7
$ws.Range("G55").Value = 14
$ws.Range("H55").Value = 48

# Row 56 -> Argentina
$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 2571
$ws.Range("D56").Value = 631
$ws.Range("E56").Value = 1825
$ws.Range("F56").Value = 121
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 115

# Row 62 (Islandia) - refreshed figures, no reorder
$ws.Range("B62").Value = 1739
$ws.Range("C62").Value = 12
$ws.Range("D62").Value = 1144
$ws.Range("E62").Value = 587
$ws.Range("F62").Value = 6

# Row 86 (Republica de Chipre) - refreshed figures, no reorder
$ws.Range("D86").Value = 77
$ws.Range("E86").Value = 646

# Row 90 (Libano) - refreshed figures, no reorder
$ws.Range("D90").Value = 86
$ws.Range("E90").Value = 556

# Row 106 (Jordania) - refreshed figures, no reorder
$ws.Range("B106").Value = 402
$ws.Range("C106").Value = 1
$ws.Range("D106").Value = 259
$ws.Range("E106").Value = 136

# Row 128 -> Jamaica
$ws.Range("A128").Value = "Jamaica"
$ws.Range("B128").Value = 143
$ws.Range("C128").Value = 18
$ws.Range("D128").Value = 21
$ws.Range("E128").Value = 117
$ws.Range("H128").Value = 5

# Row 129 -> Ruanda
$ws.Range("A129").Value = "Ruanda"
$ws.Range("B129").Value = 138
$ws.Range("C129").Value = 2
$ws.Range("D129").Value = 60
$ws.Range("E129").Value = 78
$ws.Range("F129").Value = 0
$ws.Range("H129").Value = 0

# Row 130 -> Brunei
$ws.Range("A130").Value = "Brunei"
$ws.Range("B130").Value = 136
$ws.Range("D130").Value = 108
$ws.Range("F130").Value = 2
$ws.Range("H130").Value = 1

# Row 131 -> Gibraltar
$ws.Range("A131").Value = "Gibraltar"
$ws.Range("B131").Value = 131
$ws.Range("D131").Value = 104
$ws.Range("E131").Value = 27
$ws.Range("F131").Value = 1
$ws.Range("H131").Value = 0

# Row 136 -> Aruba
$ws.Range("A136").Value = "Aruba"
$ws.Range("B136").Value = 95
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 39
$ws.Range("E136").Value = 54
$ws.Range("F136").Value = 1
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 2

# Row 137 -> Tanzania
$ws.Range("A137").Value = "Tanzania"
$ws.Range("B137").Value = 94
$ws.Range("C137").Value = 6
$ws.Range("D137").Value = 11
$ws.Range("E137").Value = 79
$ws.Range("F137").Value = 0
$ws.Range("H137").Value = 4

# Row 138 -> Monaco
$ws.Range("A138").Value = "Monaco"
$ws.Range("D138").Value = 12
$ws.Range("E138").Value = 78
$ws.Range("F138").Value = 2
$ws.Range("H138").Value = 3

# Row 153 -> Guyana
$ws.Range("A153").Value = "Guyana"
$ws.Range("D153").Value = 8
$ws.Range("E153").Value = 41
$ws.Range("F153").Value = 5
$ws.Range("H153").Value = 6

# Row 154 -> Uganda
$ws.Range("A154").Value = "Uganda"
$ws.Range("D154").Value = 20
$ws.Range("E154").Value = 35
$ws.Range("F154").Value = 0
$ws.Range("H154").Value = 0
